$d = $word.ActiveDocument

# Paragraph: "Functions of the app: " - normal paragraph (not a list item)
$p = $d.Paragraphs.Last
$null = $p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Alignment = 0
$p.Range.Text = "Functions of the app: "
$p.Range.Font.Size = 12
$p.Range.Font.SizeBi = 12

$items = @(
  "Type in a plate you want and the app pulls up ingredients you need. ",
  "Stores ingredients on a shopping list. ",
  "Have a virtual pantry to keep up what ingredients you have. ",
  "Look for plates with ingredients you have in your pantry. ",
  "Provides a calendar for users to track their meals for the week. ",
  "Users are able to search for categories of plates they desire. "
)

foreach ($text in $items) {
  $p = $d.Paragraphs.Last
  $null = $p.Range.InsertParagraphAfter()
  $p = $d.Paragraphs.Last
  $p.Alignment = 0
  $p.Range.Text = $text
  $p.Range.Font.Size = 12
  $p.Range.Font.SizeBi = 12
  $null = $p.Range.ListFormat.ApplyNumberDefault()
}

Write-Output "done"
